$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 4.382367134094238
$ws.Range("B1").Value = 4.17560338973999
$ws.Range("C1").Value = 3.65787672996521
$ws.Range("D1").Value = 1.848637819290161
$ws.Range("E1").Value = 0.8911393880844116
